$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.165.93'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '1.989.75'
$ws.Range("E3").Value = '  +5.61%  '
$c = $ws.Range("D4")
$c.Value = "'0.9983"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$c = $ws.Range("D5")
$c.Value = "'0.7997"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +68.67%  '
$c = $ws.Range("D6")
$c.Value = "'254.82"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.42%  '
$c = $ws.Range("D7")
$c.Value = "'0.9982"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '
$c = $ws.Range("D9")
$c.Value = "'28.11"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +25.99%  '
$c = $ws.Range("D10")
$c.Value = "'0.06991"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +6.85%  '
$c = $ws.Range("D11")
$c.Value = "'0.8459"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +8.95%  '
$c = $ws.Range("D12")
$c.Value = "'0.08175"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D13")
$c.Value = "'100.37"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.987.44'
$ws.Range("E14").Value = '  +5.46%  '
$c = $ws.Range("D15")
$c.Value = "'5.616"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +6.83%  '
$ws.Range("E16").Value = '  +16.71%  '
$c = $ws.Range("D17")
$c.Value = "'272.91"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -4.07%  '
$ws.Range("D18").Value = '31.154.38'
$ws.Range("E18").Value = '  +1.98%  '
$c = $ws.Range("D19")
$c.Value = "'5.869"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +9.45%  '
$c = $ws.Range("D20")
$c.Value = "'0.000007935"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.39%  '
$ws.Range("D21").Value = '2.252.10'
$ws.Range("E21").Value = '  +5.88%  '
$c = $ws.Range("D22")
$c.Value = "'0.9981"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '
$c = $ws.Range("D23")
$c.Value = "'0.9986"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '
$c = $ws.Range("D24")
$c.Value = "'7.058"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +9.38%  '
$c = $ws.Range("D25")
$c.Value = "'9.980"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +8.84%  '
$c = $ws.Range("D26")
$c.Value = "'0.1515"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +56.07%  '
$c = $ws.Range("D27")
$c.Value = "'165.73"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.69%  '
$c = $ws.Range("D28")
$c.Value = "'19.89"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +3.99%  '
$c = $ws.Range("D29")
$c.Value = "'2.342"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +22.23%  '
$ws.Range("E30").Value = '  +6.13%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D31")
$c.Value = "'4.583"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +7.64%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D32")
$c.Value = "'1.353"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.42%  '
$c = $ws.Range("D33")
$c.Value = "'4.412"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +5.26%  '
$c = $ws.Range("D34")
$c.Value = "'0.05256"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +8.41%  '
$c = $ws.Range("D35")
$c.Value = "'0.7789"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +11.59%  '
$c = $ws.Range("D36")
$c.Value = "'1.216"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +7.59%  '
$c = $ws.Range("D37")
$c.Value = "'2.758"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$c = $ws.Range("D38")
$c.Value = "'0.9972"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.28%  '
$c = $ws.Range("D39")
$c.Value = "'0.02000"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.52%  '
$c = $ws.Range("D40")
$c.Value = "'2.894"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.21%  '
$c = $ws.Range("D41")
$c.Value = "'6.644"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +5.55%  '
$c = $ws.Range("D42")
$c.Value = "'79.61"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +4.70%  '
$c = $ws.Range("D43")
$c.Value = "'0.4664"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +9.54%  '
$ws.Range("E44").Value = '  +7.04%  '
$c = $ws.Range("D45")
$c.Value = "'0.8516"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.50%  '
$c = $ws.Range("D46")
$c.Value = "'104.57"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.89%  '
$c = $ws.Range("D47")
$c.Value = "'0.9984"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$c = $ws.Range("D48")
$c.Value = "'7.678"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +9.32%  '
$c = $ws.Range("D49")
$c.Value = "'9.878"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '
$c = $ws.Range("D50")
$c.Value = "'36.82"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +4.56%  '
$c = $ws.Range("D51")
$c.Value = "'0.4292"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +8.30%  '
